# Apply updates to the "Bijgewerkte Factuur template" sheet (the active sheet):
#  - Rename the "Tarief" header (C1) to "Tarief per uur"
#  - Fill in a rate value (C2 = 10), centered to match the header column style
#  - Clear the stray "EINSCHOONMAAK AANBETALING" helper cell (F4)
#  - Clear the stray "Bij" helper cell (G16)
#  - Relabel the "Extra Kosten" rows (A10:A12) to reuse the real cost items
#  - Leave the A10:A12 range selected, matching the user's last interaction

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Tarief per uur"

$ws.Range("C2").Value = 10
$ws.Range("C2").HorizontalAlignment = -4108  # xlCenter

$ws.Range("F4").ClearContents()
$ws.Range("G16").ClearContents()

$ws.Range("A10").Value = "Schoonmaak + afvoeren afval"
$ws.Range("A11").Value = "stofzuiger"
$ws.Range("A12").Value = "Reiskosten"

$ws.Range("A10:A12").Select() | Out-Null
